$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1819908.6
$ws.Range("J17").Value = 1895700.6
$ws.Range("L17").Value = 5687101.800000001
$ws.Range("N17").Value = -5687437.800000001
$ws.Range("H40").Value = 5712.5713
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H98").Value = 72298.25
$ws.Range("I98").Value = 94231.664
$ws.Range("K98").Value = 94231.664
$ws.Range("M98").Value = -92733.664
$ws.Range("H106").Value = 6505029.5
$ws.Range("I106").Value = 7722848
$ws.Range("J106").Value = 9998.333000000001
$ws.Range("K106").Value = 7722848
$ws.Range("L106").Value = 9998.333000000001
$ws.Range("M106").Value = -7722217
$ws.Range("N106").Value = -11260.333
$ws.Range("H116").Value = 1013162.2
$ws.Range("I116").Value = 1391661.6
$ws.Range("K116").Value = 1391661.6
$ws.Range("M116").Value = -1388219.6
$ws.Range("H122").Value = 72298.25
$ws.Range("I122").Value = 94231.664
$ws.Range("K122").Value = 282694.992
$ws.Range("M122").Value = -280244.992
$ws.Range("H132").Value = 1821752.8
$ws.Range("I132").Value = 3775.1365
$ws.Range("K132").Value = 11325.4095
$ws.Range("M132").Value = -8795.4095
$ws.Range("H135").Value = 3632.077
$ws.Range("I135").Value = 3784.95
$ws.Range("J135").Value = 3122.5
$ws.Range("K135").Value = 34064.55
$ws.Range("L135").Value = 28102.5
$ws.Range("M135").Value = -31529.55
$ws.Range("N135").Value = -33172.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9418.166999999999
$ws.Range("I32").Value = 8768.759
$ws.Range("K32").Value = 8768.759
$ws.Range("M32").Value = -8481.759
$ws.Range("H43").Value = 21906.334
$ws.Range("I43").Value = 18342
$ws.Range("J43").Value = 23688.5
$ws.Range("K43").Value = 18342
$ws.Range("L43").Value = 23688.5
$ws.Range("M43").Value = -18029
$ws.Range("N43").Value = -24314.5
$ws.Range("H61").Value = 9585.645500000001
$ws.Range("I61").Value = 10598.046
$ws.Range("K61").Value = 10598.046
$ws.Range("M61").Value = -10386.046
$ws.Range("H63").Value = 1800
$ws.Range("J63").Value = 1800
$ws.Range("L63").Value = 1800
$ws.Range("N63").Value = -3172
$ws.Range("H66").Value = 1800
$ws.Range("J66").Value = 1800
$ws.Range("L66").Value = 9000
$ws.Range("N66").Value = -15864
$ws.Range("H109").Value = 75000
$ws.Range("J109").Value = 75000
$ws.Range("L109").Value = 75000
$ws.Range("N109").Value = -77774
$ws.Range("H132").Value = 1694.0294
$ws.Range("I132").Value = 1060.52
$ws.Range("J132").Value = 3453.7778
$ws.Range("K132").Value = 3181.56
$ws.Range("L132").Value = 10361.3334
$ws.Range("M132").Value = -651.5599999999999
$ws.Range("N132").Value = -15421.3334
$ws.Range("H136").Value = 9585.645500000001
$ws.Range("I136").Value = 10598.046
$ws.Range("K136").Value = 31794.138
$ws.Range("M136").Value = -29244.138

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 84849.5
$ws.Range("J35").Value = 84849.5
$ws.Range("L35").Value = 84849.5
$ws.Range("N35").Value = -85469.5
$ws.Range("H105").Value = 64453.11
$ws.Range("I105").Value = 101843.27
$ws.Range("K105").Value = 101843.27
$ws.Range("M105").Value = -100096.27
$ws.Range("H134").Value = 2244
$ws.Range("I134").Value = 1317.45
$ws.Range("K134").Value = 3952.35
$ws.Range("M134").Value = -1417.35

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 406.25
$ws.Range("I22").Value = 395
$ws.Range("K22").Value = 395
$ws.Range("M22").Value = -45
$ws.Range("H31").Value = 6101.975
$ws.Range("I31").Value = 6227.121
$ws.Range("K31").Value = 6227.121
$ws.Range("M31").Value = -5932.121
$ws.Range("H34").Value = 6101.975
$ws.Range("I34").Value = 6227.121
$ws.Range("K34").Value = 6227.121
$ws.Range("M34").Value = -6025.121
$ws.Range("H58").Value = 2522.9119
$ws.Range("J58").Value = 3383.8572
$ws.Range("L58").Value = 3383.8572
$ws.Range("N58").Value = -3789.8572
$ws.Range("H136").Value = 2522.9119
$ws.Range("J136").Value = 3383.8572
$ws.Range("L136").Value = 10151.5716
$ws.Range("N136").Value = -15251.5716
$ws.Range("H141").Value = 182792.7
$ws.Range("I141").Value = 71763.664
$ws.Range("J141").Value = 198654
$ws.Range("K141").Value = 71763.664
$ws.Range("L141").Value = 198654
$ws.Range("M141").Value = -66583.664
$ws.Range("N141").Value = -209014

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -38744
$ws.Range("H122").Value = 11555.529
$ws.Range("I122").Value = 8057.231
$ws.Range("J122").Value = 22925
$ws.Range("K122").Value = 24171.693
$ws.Range("L122").Value = 68775
$ws.Range("M122").Value = -21721.693
$ws.Range("N122").Value = -73675
$ws.Range("H126").Value = 7170.8066
$ws.Range("I126").Value = 14298.5
$ws.Range("J126").Value = 3776.6667
$ws.Range("K126").Value = 42895.5
$ws.Range("L126").Value = 11330.0001
$ws.Range("M126").Value = -40425.5
$ws.Range("N126").Value = -16270.0001
$ws.Range("H132").Value = 1987.4524
$ws.Range("I132").Value = 1853.5135
$ws.Range("J132").Value = 2978.6
$ws.Range("K132").Value = 5560.5405
$ws.Range("L132").Value = 8935.799999999999
$ws.Range("M132").Value = -3030.5405
$ws.Range("N132").Value = -13995.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 11395
$ws.Range("I20").Value = 11395
$ws.Range("K20").Value = 11395
$ws.Range("M20").Value = -11169
$ws.Range("H22").Value = 666.75
$ws.Range("I22").Value = 666.75
$ws.Range("K22").Value = 666.75
$ws.Range("M22").Value = -371.75
$ws.Range("H27").Value = 666.75
$ws.Range("I27").Value = 666.75
$ws.Range("K27").Value = 666.75
$ws.Range("M27").Value = -559.75
$ws.Range("H40").Value = 22001.965
$ws.Range("I40").Value = 28703.555
$ws.Range("K40").Value = 28703.555
$ws.Range("M40").Value = -28567.555
$ws.Range("H46").Value = 2415.5217
$ws.Range("I46").Value = 710.125
$ws.Range("K46").Value = 710.125
$ws.Range("M46").Value = -522.125
$ws.Range("H132").Value = 575467.75
$ws.Range("I132").Value = 746470.9399999999
$ws.Range("K132").Value = 2239412.82
$ws.Range("M132").Value = -2236882.82

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 15218.851
$ws.Range("I122").Value = 2414.1765
$ws.Range("K122").Value = 7242.529500000001
$ws.Range("M122").Value = -4792.529500000001
